$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C10 changes from 18 to 1 (numeric value)
$ws.Range("C10").Value = 1
